# Weekly update: insert a new "latest" record for Espárragos (Macroferia
# Regional de Talca) at row 19, pushing the existing rows 19-23 down to
# rows 20-24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 19; this shifts the old
# rows 19-23 down to 20-24 (mirrors the pattern in the target diff).
$ws.Rows("19:19").Insert()

# Populate the newly inserted row 19 with this week's data. Columns
# A, B, C, E, F, G, H, I, O, Q, R are constant for every record in this
# subset (same market / region / product / origin), so just repeat them.
$ws.Range("A19").Value = 5
$ws.Range("B19").Value = "Macroferia Regional de Talca"
$ws.Range("C19").Value = "Maule"
$ws.Range("D19").Value = 44474
$ws.Range("E19").Value = 7
$ws.Range("F19").Value = 300000000
$ws.Range("G19").Value = "Espárragos"
$ws.Range("H19").Value = "Verde"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 5000
$ws.Range("K19").Value = 1200
$ws.Range("L19").Value = 1200
$ws.Range("M19").Value = 1200
$ws.Range("N19").Value = "$/kilo"
$ws.Range("O19").Value = "Provincia de Linares"
$ws.Range("P19").Value = 1200
$ws.Range("Q19").Value = 1
$ws.Range("R19").Value = "Hortaliza"
